# US 3.3 commit files
# Update the BGDPbES (BAU Guaranteed Dispatch Perc by Elec Source) workbook:
#  - add a header label + formatting to A1 on the BGDPbES sheet
#  - rename the "coal to gas" fuel row to "lignite"
#  - turn on guaranteed dispatch (set to 1) for nuclear and hydro
#  - append three new fuel rows: crude oil, heavy or residual fuel oil,
#    municipal solid waste
#  - restore the "About" sheet as the active tab

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGDPbES")

# --- New header cell on row 1 ------------------------------------------
$ws.Range("A1").Value = "BAU Guaranteed Dispatch (dimensionless)"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 45

# --- Rename the "coal to gas" row to "lignite" --------------------------
$ws.Range("A13").Value = "lignite"

# --- Turn on guaranteed dispatch for nuclear (row 4) and hydro (row 5) --
$ws.Range("B4").Value = 1
$ws.Range("B5").Value = 1

# --- Append three new fuel-source rows (plain values, no formulas) ------
$ws.Range("A15").Value = "crude oil"
$ws.Range("B15:AK15").Value = 0

$ws.Range("A16").Value = "heavy or residual fuel oil"
$ws.Range("B16:AK16").Value = 0

$ws.Range("A17").Value = "municipal solid waste"
$ws.Range("B17:AK17").Value = 0

# --- Page setup tweak on the BGDPbES sheet -------------------------------
$ws.PageSetup.Orientation = 1

# --- Make "About" the active sheet/tab again -----------------------------
$ws.Range("A1").Select()
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$wsAbout.Range("A1").Select()
